# Update "想去人数" (want-to-go count) figures and mark two tickets as
# "不可售" (not sellable) across the relevant worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 17
$ws1.Range("F4").Value = 1071
$ws1.Range("F5").Value = 19749
$ws1.Range("F7").Value = 2269
$ws1.Range("F11").Value = 694
$ws1.Range("F12").Value = 239
$ws1.Range("F13").Value = 253
$ws1.Range("F15").Value = 372
$ws1.Range("F17").Value = 264
$ws1.Range("F19").Value = 195

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F7").Value = 289
$ws2.Range("F15").Value = 77
$ws2.Range("F20").Value = 19

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6005
$ws3.Range("F4").Value = 593

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6005
$ws4.Range("F4").Value = 593
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F6").Value = 17
$ws4.Range("F9").Value = 1071
$ws4.Range("F10").Value = 19749
$ws4.Range("F15").Value = 289
$ws4.Range("F16").Value = 2269
$ws4.Range("F21").Value = 694
$ws4.Range("F22").Value = 239
$ws4.Range("F23").Value = 253
$ws4.Range("F28").Value = 372
$ws4.Range("F31").Value = 264
$ws4.Range("F35").Value = 195
$ws4.Range("F36").Value = 77
$ws4.Range("F43").Value = 19
